$wb = $excel.ActiveWorkbook

# "Repayment Schedule" is the 3rd worksheet (sheet3.xml) - insert a new blank
# column before column N (shifting old N "Late", O "heading", P "Outstanding"
# one column to the right), reflecting the new Variable Instalments layout.
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Columns("N:N").Insert()

# Make "Repayment Schedule" the active/selected sheet with the given selection,
# which also clears the "Output" sheet's tabSelected flag.
$ws.Activate()
$ws.Range("J25").Select()
